$wb = $excel.ActiveWorkbook

# Sheet1 (Planilha1): only the selection changes (view state) - D3 -> B2:B7
$ws1 = $wb.Worksheets.Item("Planilha1")
$ws1.Range("B2:B7").Select()

# Sheet2 (Planilha2): C2 date value changes 44712 -> 44713 (2022-05-31 -> 2022-06-01)
# and a new row 8 with B8 cell (empty, percentage style like the rest of column B)
$ws2 = $wb.Worksheets.Item("Planilha2")
$ws2.Range("C2").Value = 44713
$ws2.Range("B8").NumberFormat = $ws2.Range("B7").NumberFormat
$ws2.Range("D5").Select()

# Sheet3 (regras): the note text moves from B2 to A1, with updated wording
$ws3 = $wb.Worksheets.Item("regras")
$ws3.Range("B2").ClearContents()
$ws3.Range("A1").Value = "A data inicio deve ser no dia seguinte da data fim anterior"
$ws3.Activate()
$ws3.Range("A1").Select()
